$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# J116/K116 new values
$ws.Cells.Item(116,10).Value = 35610
$ws.Cells.Item(116,11).Value = 35604

# B122/B123 updated values
$ws.Cells.Item(122,2).Value = 37222
$ws.Cells.Item(123,2).Value = 37448

# Row 124 (fill in A/B, C124 already exists, add D formula)
$ws.Cells.Item(124,1).Value = "Enter Pipe"
$ws.Cells.Item(124,2).Value = 37664
$ws.Cells.Item(124,4).Formula = "=IF(B124 >  0,C124-B124, 0)"

# Row 125
$ws.Cells.Item(125,1).Value = "Checkpoint Rail 12582912"
$ws.Cells.Item(125,2).Value = 37919
$ws.Cells.Item(125,3).Value = 44527
$ws.Cells.Item(125,4).Formula = "=IF(B125 >  0,C125-B125, 0)"

# Row 126
$ws.Cells.Item(126,1).Value = "Checkpoint 460"
$ws.Cells.Item(126,2).Value = 39165
$ws.Cells.Item(126,3).Value = 45775
$ws.Cells.Item(126,4).Formula = "=IF(B126 >  0,C126-B126, 0)"

# Row 127
$ws.Cells.Item(127,1).Value = "Cehckpiont 1449"
$ws.Cells.Item(127,2).Value = 39722
$ws.Cells.Item(127,3).Value = 46333
$ws.Cells.Item(127,4).Formula = "=IF(B127 >  0,C127-B127, 0)"

# Row 128
$ws.Cells.Item(128,1).Value = "Enter Pipe"
$ws.Cells.Item(128,2).Value = 41771
$ws.Cells.Item(128,3).Value = 48401
$ws.Cells.Item(128,4).Formula = "=IF(B128 >  0,C128-B128, 0)"

# Row 129
$ws.Cells.Item(129,1).Value = "Get flag"
$ws.Cells.Item(129,2).Value = 42012
$ws.Cells.Item(129,3).Value = 48651
$ws.Cells.Item(129,4).Formula = "=IF(B129 >  0,C129-B129, 0)"

# Row 130
$ws.Cells.Item(130,1).Value = "End level"
$ws.Cells.Item(130,2).Value = 42526
$ws.Cells.Item(130,3).Value = 49165
$ws.Cells.Item(130,4).Formula = "=IF(B130 >  0,C130-B130, 0)"

# Row 131
$ws.Cells.Item(131,1).Value = "Enter 8-4"
$ws.Cells.Item(131,2).Value = 42919
$ws.Cells.Item(131,3).Value = 49956
$ws.Cells.Item(131,4).Formula = "=IF(B131 >  0,C131-B131, 0)"

# Row 132
$ws.Cells.Item(132,1).Value = "1st Move"
$ws.Cells.Item(132,2).Value = 43147
$ws.Cells.Item(132,3).Value = 50203
$ws.Cells.Item(132,4).Formula = "=IF(B132 >  0,C132-B132, 0)"

$ws.Range("B133").Select()
